$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 933.3333
$ws.Range("I19").Value = 800
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 800
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = -625
$ws.Range("N19").Value = -1350
$ws.Range("H62").Value = 5633
$ws.Range("I62").Value = 3942.3572
$ws.Range("J62").Value = 9014.286
$ws.Range("K62").Value = 3942.3572
$ws.Range("L62").Value = 9014.286
$ws.Range("M62").Value = -3318.3572
$ws.Range("N62").Value = -10262.286
$ws.Range("H65").Value = 5633
$ws.Range("I65").Value = 3942.3572
$ws.Range("J65").Value = 9014.286
$ws.Range("K65").Value = 19711.786
$ws.Range("L65").Value = 45071.43
$ws.Range("M65").Value = -16591.786
$ws.Range("N65").Value = -51311.43
$ws.Range("H129").Value = 1000
$ws.Range("I129").Value = 390
$ws.Range("J129").Value = 1110.909
$ws.Range("K129").Value = 1170
$ws.Range("L129").Value = 3332.727
$ws.Range("M129").Value = 3830
$ws.Range("N129").Value = -13332.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2679.2068
$ws.Range("I32").Value = 2722.55
$ws.Range("J32").Value = 2582.889
$ws.Range("K32").Value = 2722.55
$ws.Range("L32").Value = 2582.889
$ws.Range("M32").Value = -2435.55
$ws.Range("N32").Value = -3156.889
$ws.Range("H61").Value = 2020.2858
$ws.Range("I61").Value = 1616.75
$ws.Range("K61").Value = 1616.75
$ws.Range("M61").Value = -1404.75
$ws.Range("H74").Value = 1047.8572
$ws.Range("I74").Value = 1203.0769
$ws.Range("J74").Value = 913.3333
$ws.Range("K74").Value = 1203.0769
$ws.Range("L74").Value = 913.3333
$ws.Range("M74").Value = -329.0769
$ws.Range("N74").Value = -2661.3333
$ws.Range("H77").Value = 1047.8572
$ws.Range("I77").Value = 1203.0769
$ws.Range("J77").Value = 913.3333
$ws.Range("K77").Value = 6015.3845
$ws.Range("L77").Value = 4566.6665
$ws.Range("M77").Value = -1647.3845
$ws.Range("N77").Value = -13302.6665
$ws.Range("H136").Value = 2020.2858
$ws.Range("I136").Value = 1616.75
$ws.Range("K136").Value = 4850.25
$ws.Range("M136").Value = -2300.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 635.875
$ws.Range("I94").Value = 611.6
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 611.6
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -160.6
$ws.Range("N94").Value = -1902
$ws.Range("H132").Value = 16263.333
$ws.Range("J132").Value = 16263.333
$ws.Range("L132").Value = 16263.333
$ws.Range("N132").Value = -26383.333
$ws.Range("H134").Value = 2386.3572
$ws.Range("I134").Value = 1475.6786
$ws.Range("J134").Value = 4207.7144
$ws.Range("K134").Value = 4427.0358
$ws.Range("L134").Value = 12623.1432
$ws.Range("M134").Value = -1892.0358
$ws.Range("N134").Value = -17693.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1685.1034
$ws.Range("I58").Value = 996
$ws.Range("J58").Value = 2245
$ws.Range("K58").Value = 996
$ws.Range("L58").Value = 2245
$ws.Range("M58").Value = -793
$ws.Range("N58").Value = -2651
$ws.Range("H60").Value = 13550
$ws.Range("I60").Value = 4540
$ws.Range("J60").Value = 28566.666
$ws.Range("K60").Value = 4540
$ws.Range("L60").Value = 28566.666
$ws.Range("M60").Value = -4029
$ws.Range("N60").Value = -29588.666
$ws.Range("H132").Value = 2406.543
$ws.Range("I132").Value = 1325
$ws.Range("K132").Value = 3975
$ws.Range("M132").Value = -1445
$ws.Range("H136").Value = 1685.1034
$ws.Range("I136").Value = 996
$ws.Range("J136").Value = 2245
$ws.Range("K136").Value = 2988
$ws.Range("L136").Value = 6735
$ws.Range("M136").Value = -438
$ws.Range("N136").Value = -11835

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 4000
$ws.Range("J17").Value = 4000
$ws.Range("L17").Value = 12000
$ws.Range("N17").Value = -12338
$ws.Range("H108").Value = 3518.6365
$ws.Range("I108").Value = 2967.2222
$ws.Range("K108").Value = 8901.6666
$ws.Range("M108").Value = -6021.6666
$ws.Range("H113").Value = 401.02856
$ws.Range("J113").Value = 419.08694
$ws.Range("L113").Value = 1257.26082
$ws.Range("N113").Value = -5597.26082
$ws.Range("H115").Value = 1649.75
$ws.Range("I115").Value = 1199.6666
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 3598.9998
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -2423.9998
$ws.Range("N115").Value = -11350
$ws.Range("H118").Value = 1369.75
$ws.Range("I118").Value = 1165.4286
$ws.Range("J118").Value = 2800
$ws.Range("K118").Value = 3496.2858
$ws.Range("L118").Value = 8400
$ws.Range("M118").Value = -2253.2858
$ws.Range("N118").Value = -10886
$ws.Range("H121").Value = 165
$ws.Range("I121").Value = 165
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 495
$ws.Range("L121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = 815
$ws.Range("H131").Value = 2491.6748
$ws.Range("J131").Value = 2771.7397
$ws.Range("L131").Value = 8315.2191
$ws.Range("N131").Value = -18395.2191

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1835.5
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 2171
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2171
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6511
$ws.Range("H126").Value = 2067.238
$ws.Range("I126").Value = 985.3333
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 2955.9999
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -485.9998999999998
$ws.Range("N126").Value = -12440
$ws.Range("H132").Value = 2570.3953
$ws.Range("I132").Value = 2166.484
$ws.Range("J132").Value = 3613.8333
$ws.Range("K132").Value = 6499.451999999999
$ws.Range("L132").Value = 10841.4999
$ws.Range("M132").Value = -3969.451999999999
$ws.Range("N132").Value = -15901.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3167.76
$ws.Range("I7").Value = 2313.4285
$ws.Range("K7").Value = 2313.4285
$ws.Range("M7").Value = -2201.4285
$ws.Range("H45").Value = 14515
$ws.Range("I45").Value = 10041
$ws.Range("J45").Value = 18989
$ws.Range("K45").Value = 10041
$ws.Range("L45").Value = 18989
$ws.Range("M45").Value = -9634
$ws.Range("N45").Value = -19803
$ws.Range("H126").Value = 3167.76
$ws.Range("I126").Value = 2313.4285
$ws.Range("K126").Value = 6940.2855
$ws.Range("M126").Value = -4470.2855
$ws.Range("H132").Value = 3125.4517
$ws.Range("I132").Value = 2457.3462
$ws.Range("J132").Value = 6599.6
$ws.Range("K132").Value = 7372.0386
$ws.Range("L132").Value = 19798.8
$ws.Range("M132").Value = -4842.0386
$ws.Range("N132").Value = -24858.8
$ws.Range("H136").Value = 7287
$ws.Range("I136").Value = 3212.5
$ws.Range("K136").Value = 9637.5
$ws.Range("M136").Value = -7087.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 20897412
$ws.Range("I136").Value = 41792396
$ws.Range("J136").Value = 2427.125
$ws.Range("K136").Value = 125377188
$ws.Range("L136").Value = 7281.375
$ws.Range("M136").Value = -125374638
$ws.Range("N136").Value = -12381.375
